$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column widths (A=26, B=14, C=6 characters) ---------------------------
# The engine quantizes ColumnWidth to pixel steps; empirically the first
# ColumnWidth value that rounds to an exact integer N (character width) in
# the saved <col width="N"/> is N - 11/12. Using that offset reproduces the
# exact target widths of 26 / 14 / 6.
$ws.Columns.Item(1).ColumnWidth = 26 - 11/12
$ws.Columns.Item(2).ColumnWidth = 14 - 11/12
$ws.Columns.Item(3).ColumnWidth = 6 - 11/12

# --- Data fixes -------------------------------------------------------------
# Column B: drop the leading ". " before the registration number, keeping a
# single leading space (". GT24BCAR001" -> " GT24BCAR001").
# Column C: convert the SGPA values from text to real numbers, with the
# placeholder "--" becoming 0.
$regNos = @(
    "GT24BCAR001",
    "GT24BCAR002",
    "GT24BCAR003",
    "GT24BCAR006",
    "GT24BCAR007",
    "GT24BCAR008",
    "GT24BCAR009",
    "GT24BCAR011",
    "GT24BCAR013",
    "GT24BCAR014",
    "GT24BCAR015",
    "GT24BCAR016",
    "GT24BCAR017",
    "GT24BCAR018",
    "GT24BCAR019",
    "GT24BCAR020",
    "GT24BCAR021",
    "GT24BCAR023",
    "GT24BCAR024",
    "GT24BCAR025",
    "GT24BCAR026",
    "GT24BCAR027",
    "GT24BCAR028",
    "GT24BCAR029",
    "GT24BCAR030",
    "GT24BCAR031",
    "GT24BCAR032",
    "GT24BCAR034",
    "GT24BCAR004",
    "GT24BCAR022",
    "GT24BCAR033",
    "GT24BCAR035"
)

$sgpa = @(
    8.52,
    6.43,
    0,
    8.10,
    6.90,
    0,
    7.71,
    7.00,
    6.19,
    0,
    8.05,
    6.76,
    7.95,
    7.19,
    7.95,
    8.24,
    6.57,
    6.33,
    6.90,
    0,
    6.90,
    7.71,
    8.62,
    6.14,
    0,
    6.05,
    7.19,
    7.00,
    7.24,
    6.76,
    7.33,
    8.52
)

for ($i = 0; $i -lt $regNos.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 2).Value = " " + $regNos[$i]
    $ws.Cells.Item($row, 3).Value = $sgpa[$i]
}
